# Rename the "displayName" header (column B, row 1) to "display_name".
# The underlying data for every other cell is unchanged - only the
# shared-string table is affected (the old "displayName" entry is
# replaced by a new "display_name" entry), which Excel re-indexes on
# save automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "display_name"

# Match the author's final selection/cursor position: cell B1 selected
# (instead of the prior D22), with the sheet scrolled back to the top
# (topLeftCell reset from A6).
$ws.Activate()
$ws.Range("B1").Select()
